# Daily attendance processing - 2026-01-04 19:30:54
# Reorders the comma-separated "Recorded By" entries in column G so that
# the displayed author order matches the freshly re-generated report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "dnasr281@gmail.com, System" = "System, dnasr281@gmail.com"
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
